$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Range("H8").Value = 91.75
$ws.Range("I8").Value = 106.76923
$ws.Range("J8").Value = 26.666666
$ws.Range("K8").Value = 320.30769
$ws.Range("L8").Value = 79.99999800000001
$ws.Range("M8").Value = -181.30769
$ws.Range("N8").Value = -357.999998
# row 19
$ws.Range("H19").Value = 2229
$ws.Range("I19").Value = 2049.6667
$ws.Range("J19").Value = 2498
$ws.Range("K19").Value = 2049.6667
$ws.Range("L19").Value = 2498
$ws.Range("M19").Value = -1874.6667
$ws.Range("N19").Value = -2848
# row 33
$ws.Range("H33").Value = 1064.1666
$ws.Range("I33").Value = 75.8
$ws.Range("J33").Value = 6006
$ws.Range("K33").Value = 75.8
$ws.Range("L33").Value = 6006
$ws.Range("M33").Value = 153.2
$ws.Range("N33").Value = -6464
# row 40
$ws.Range("H40").Value = 1408.3636
$ws.Range("I40").Value = 1261
$ws.Range("J40").Value = 1492.5714
$ws.Range("K40").Value = 1261
$ws.Range("L40").Value = 1492.5714
$ws.Range("M40").Value = -1086
$ws.Range("N40").Value = -1842.5714
# row 53
$ws.Range("H53").Value = 326.22223
$ws.Range("I53").Value = 1010
$ws.Range("J53").Value = 130.85715
$ws.Range("K53").Value = 1010
$ws.Range("L53").Value = 130.85715
$ws.Range("M53").Value = -373
$ws.Range("N53").Value = -1404.85715
# row 100
$ws.Range("H100").Value = 2798.7
$ws.Range("I100").Value = 2869.5715
$ws.Range("J100").Value = 2633.3333
$ws.Range("K100").Value = 2869.5715
$ws.Range("L100").Value = 2633.3333
$ws.Range("M100").Value = -2328.5715
$ws.Range("N100").Value = -3715.3333
# row 112
$ws.Range("H112").Value = 2219.2
$ws.Range("I112").Value = 1365.3334
$ws.Range("J112").Value = 3500
$ws.Range("K112").Value = 4096.0002
$ws.Range("L112").Value = 10500
$ws.Range("M112").Value = -2988.0002
$ws.Range("N112").Value = -12716
# row 137
$ws.Range("H137").Value = 1873.7084
$ws.Range("I137").Value = 1597.55
$ws.Range("J137").Value = 3254.5
$ws.Range("K137").Value = 4792.65
$ws.Range("L137").Value = 9763.5
$ws.Range("M137").Value = -2242.65
$ws.Range("N137").Value = -14863.5
# row 138
$ws.Range("H138").Value = 2424.2666
$ws.Range("I138").Value = 1394.4
$ws.Range("J138").Value = 2939.2
$ws.Range("K138").Value = 4183.200000000001
$ws.Range("L138").Value = 8817.599999999999
$ws.Range("M138").Value = 956.7999999999993
$ws.Range("N138").Value = -19097.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 117.111115
$ws.Range("I5").Value = 113
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 113
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -1
$ws.Range("N5").Value = -374
# row 32
$ws.Range("H32").Value = 1622.1892
$ws.Range("I32").Value = 1500.5834
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 1500.5834
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -1213.5834
$ws.Range("N32").Value = -6574
# row 55
$ws.Range("H55").Value = 28333.334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 28333.334
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 28333.334
$ws.Range("N55").Value = -28963.334
# row 61
$ws.Range("H61").Value = 5737.3335
$ws.Range("I61").Value = 5737.3335
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5737.3335
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5525.3335
# row 74
$ws.Range("H74").Value = 1399.5
$ws.Range("I74").Value = 1399.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1399.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -525.5
# row 77
$ws.Range("H77").Value = 1399.5
$ws.Range("I77").Value = 1399.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6997.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2629.5
# row 122
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -8500
# row 132
$ws.Range("H132").Value = 2456.75
$ws.Range("I132").Value = 2456.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7370.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4840.25
$ws.Range("N132").ClearContents()
# row 136
$ws.Range("H136").Value = 5737.3335
$ws.Range("I136").Value = 5737.3335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 17212.0005
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14662.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 117.111115
$ws.Range("I4").Value = 113
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 113
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = -380
# row 19
$ws.Range("H19").Value = 50010
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 50010
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 50010
$ws.Range("N19").Value = -50356
# row 107
$ws.Range("H107").Value = 2440
$ws.Range("I107").Value = 2566.6667
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 2566.6667
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = -646.6667000000002
$ws.Range("N107").Value = -6090

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 3009.2666
$ws.Range("I58").Value = 2812.7778
$ws.Range("J58").Value = 3304
$ws.Range("K58").Value = 2812.7778
$ws.Range("L58").Value = 3304
$ws.Range("M58").Value = -2609.7778
$ws.Range("N58").Value = -3710
# row 86
$ws.Range("H86").Value = 6692.75
$ws.Range("I86").Value = 6403.3335
$ws.Range("J86").Value = 7064.857
$ws.Range("K86").Value = 6403.3335
$ws.Range("L86").Value = 7064.857
$ws.Range("M86").Value = -5280.3335
$ws.Range("N86").Value = -9310.857
# row 89
$ws.Range("H89").Value = 6692.75
$ws.Range("I89").Value = 6403.3335
$ws.Range("J89").Value = 7064.857
$ws.Range("K89").Value = 32016.6675
$ws.Range("L89").Value = 35324.285
$ws.Range("M89").Value = -26400.6675
$ws.Range("N89").Value = -46556.285
# row 105
$ws.Range("H105").Value = 1249.6666
$ws.Range("I105").Value = 1249.6666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1249.6666
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 497.3334
# row 132
$ws.Range("H132").Value = 2421.5557
$ws.Range("I132").Value = 2421.5557
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7264.6671
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4734.6671
# row 136
$ws.Range("H136").Value = 3009.2666
$ws.Range("I136").Value = 2812.7778
$ws.Range("J136").Value = 3304
$ws.Range("K136").Value = 8438.3334
$ws.Range("L136").Value = 9912
$ws.Range("M136").Value = -5888.3334
$ws.Range("N136").Value = -15012
# row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 29
$ws.Range("H29").Value = 218.25
$ws.Range("I29").Value = 207.66667
$ws.Range("J29").Value = 250
$ws.Range("K29").Value = 623.00001
$ws.Range("L29").Value = 750
$ws.Range("M29").Value = -346.00001
$ws.Range("N29").Value = -1304
# row 38
$ws.Range("H38").Value = 163.2
$ws.Range("I38").Value = 153.28572
$ws.Range("J38").Value = 186.33333
$ws.Range("K38").Value = 459.85716
$ws.Range("L38").Value = 558.99999
$ws.Range("M38").Value = -112.85716
$ws.Range("N38").Value = -1252.99999
# row 46
$ws.Range("H46").Value = 2993.75
$ws.Range("I46").Value = 2475
$ws.Range("J46").Value = 3166.6667
$ws.Range("K46").Value = 7425
$ws.Range("L46").Value = 9500.000100000001
$ws.Range("M46").Value = -7334
$ws.Range("N46").Value = -9682.000100000001
# row 92
$ws.Range("H92").Value = 413.25
$ws.Range("I92").Value = 411.33334
$ws.Range("J92").Value = 419
$ws.Range("K92").Value = 1234.00002
$ws.Range("L92").Value = 1257
$ws.Range("M92").Value = 13.99998000000005
$ws.Range("N92").Value = -3753
# row 107
$ws.Range("H107").Value = 716.5
$ws.Range("I107").Value = 650.1667
$ws.Range("J107").Value = 782.8333
$ws.Range("K107").Value = 1950.5001
$ws.Range("L107").Value = 2348.4999
$ws.Range("M107").Value = -30.50009999999997
$ws.Range("N107").Value = -6188.4999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 23214.285
$ws.Range("I57").Value = 12000
$ws.Range("J57").Value = 25083.334
$ws.Range("K57").Value = 12000
$ws.Range("L57").Value = 25083.334
$ws.Range("M57").Value = -11180
$ws.Range("N57").Value = -26723.334
# row 97
$ws.Range("H97").Value = 564.1875
$ws.Range("I97").Value = 581.9286
$ws.Range("J97").Value = 440
$ws.Range("K97").Value = 581.9286
$ws.Range("L97").Value = 440
$ws.Range("M97").Value = -85.92859999999996
$ws.Range("N97").Value = -1432
# row 132
$ws.Range("H132").Value = 2544.4119
$ws.Range("I132").Value = 2544.4119
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7633.2357
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5103.2357

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 19794.105
$ws.Range("I7").Value = 21622.883
$ws.Range("J7").Value = 4249.5
$ws.Range("K7").Value = 21622.883
$ws.Range("L7").Value = 4249.5
$ws.Range("M7").Value = -21510.883
$ws.Range("N7").Value = -4473.5
# row 22
$ws.Range("H22").Value = 1161.5555
$ws.Range("I22").Value = 699.3333
$ws.Range("J22").Value = 1392.6666
$ws.Range("K22").Value = 699.3333
$ws.Range("L22").Value = 1392.6666
$ws.Range("M22").Value = -404.3333
$ws.Range("N22").Value = -1982.6666
# row 27
$ws.Range("H27").Value = 1161.5555
$ws.Range("I27").Value = 699.3333
$ws.Range("J27").Value = 1392.6666
$ws.Range("K27").Value = 699.3333
$ws.Range("L27").Value = 1392.6666
$ws.Range("M27").Value = -592.3333
$ws.Range("N27").Value = -1606.6666
# row 40
$ws.Range("H40").Value = 1742.8889
$ws.Range("I40").Value = 1742.8889
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1742.8889
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1606.8889
# row 93
$ws.Range("H93").Value = 1326.875
$ws.Range("I93").Value = 1230.7142
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1230.7142
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 17.28580000000011
$ws.Range("N93").Value = -4496
# row 100
$ws.Range("H100").Value = 4999.8335
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4999.8335
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4999.8335
$ws.Range("N100").Value = -6081.8335
# row 126
$ws.Range("H126").Value = 19794.105
$ws.Range("I126").Value = 21622.883
$ws.Range("J126").Value = 4249.5
$ws.Range("K126").Value = 64868.649
$ws.Range("L126").Value = 12748.5
$ws.Range("M126").Value = -62398.649
$ws.Range("N126").Value = -17688.5
# row 132
$ws.Range("H132").Value = 17421.428
$ws.Range("I132").Value = 11000
$ws.Range("J132").Value = 19990
$ws.Range("K132").Value = 33000
$ws.Range("L132").Value = 59970
$ws.Range("M132").Value = -30470
$ws.Range("N132").Value = -65030

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 3920.8
$ws.Range("I132").Value = 3920.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11762.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9232.400000000001

Write-Output "edits applied"